$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @(
    "giặt sấy",
    "nhà thuốc tư nhân",
    "nhà thuốc tây",
    "nhà thuốc",
    "vựa khô đường đậu gia vị ",
    "tã sữa tổng hợp",
    "lan ,mỹ phẩm",
    "chuyên bán sỉ lẻ nước ngọt",
    "cưa hàng gia dụng",
    "phụ liệu tóc nail",
    "nhận sửa quần áo",
    "chuyên cung cấp các loại sữa",
    "tạp hoa gia dụng",
    "xe gắn máy",
    "ehome"
)

$startRow = 258
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}

$lastRow = $startRow + $newValues.Count - 1
$excel.ActiveWindow.ScrollRow = 260
$ws.Range("A" + $lastRow).Select()
